{"js": "// Update the benchmark stats table: the document is a single-column table\n// where each row holds one stat value. Several rows get their value text\n// swapped, and the final three rows (which previously held a whole\n// tab-separated stats line crammed into one run) collapse down to just\n// their leading count value.\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"276\",\n  5: \"0.00048\",\n  6: \"0.00007\",\n  7: \"0.00003\",\n  8: \"0.00011\",\n  9: \"0.00014\",\n  10: \"0.00014\",\n  11: \"0.01709\",\n  43: \"100\",\n  44: \"0.02\",\n  45: \"787\",\n};\n\nfor (const [rowIndex, newValue] of Object.entries(updates)) {\n  table.getCell(Number(rowIndex), 0).value = newValue;\n}\n\nawait context.sync();\n", "ps1": "# Update the benchmark stats table: the document is a single-column table\n# where each row holds one stat value. Several rows get their value text\n# swapped, and the final three rows (which previously held a whole\n# tab-separated stats line crammed into one run) collapse down to just\n# their leading count value.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"276\"\n    6  = \"0.00048\"\n    7  = \"0.00007\"\n    8  = \"0.00003\"\n    9  = \"0.00011\"\n    10 = \"0.00014\"\n    11 = \"0.00014\"\n    12 = \"0.01709\"\n    44 = \"100\"\n    45 = \"0.02\"\n    46 = \"787\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $t.Cell($rowIndex, 1).Range.Text = $updates[$rowIndex]\n}\n"}
